# "Generate Report for handback"
#
# The localization-status report is regenerated after a handback event:
#   - the two tracked files (2eff61ea...md and 78740cd6...md) move from
#     "Not yet handed off" to "Handed back" status on every language sheet
#   - each language sheet gains two new columns of data for row 2 & 3:
#       E = Latest Target File    (same file/link as column A)
#       F = Latest Handback File  (same file/link as column C)
#       G = Latest Handback DateTime (timestamp of the handback)
#
$wb = $excel.ActiveWorkbook

# Hyperlink colour used by the workbook's existing "HyperLink" cell style
# (FF6495ED == RGB(0x64,0x95,0xED) expressed as a VBA/OLE BGR long)
$linkColor = 15570276

function Set-HandbackRow {
    param(
        $ws,
        [string]$statusCell,
        [string]$targetCell,
        [string]$targetDisplay,
        [string]$targetUrl,
        [string]$handbackCell,
        [string]$handbackDisplay,
        [string]$handbackUrl,
        [string]$dateCell,
        [string]$dateValue
    )

    # Status -> "Handed back"
    $ws.Range($statusCell).Value = "Handed back"

    # Latest Target File (mirrors column A's hyperlink/display)
    $ws.Range($targetCell).Value = $targetDisplay
    $ws.Hyperlinks.Add($ws.Range($targetCell), $targetUrl, "", "", $targetDisplay)
    $ws.Range($targetCell).Style = "HyperLink"
    $ws.Range($targetCell).Font.Underline = 2
    $ws.Range($targetCell).Font.Color = $linkColor

    # Latest Handback File (mirrors column C's hyperlink/display)
    $ws.Range($handbackCell).Value = $handbackDisplay
    $ws.Hyperlinks.Add($ws.Range($handbackCell), $handbackUrl, "", "", $handbackDisplay)
    $ws.Range($handbackCell).Style = "HyperLink"
    $ws.Range($handbackCell).Font.Underline = 2
    $ws.Range($handbackCell).Font.Color = $linkColor

    # Latest Handback DateTime
    $ws.Range($dateCell).Value = $dateValue
}

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

Set-HandbackRow $wsZh "B2" "E2" "2eff61ea-3a49-4af9-a49a-115643b7abf4.md" `
    "https://github.com/OpenLocalizationTest/oltest/blob/7d4cba76265f0c0f490220592f21e3eec828d83f/e2e/2eff61ea-3a49-4af9-a49a-115643b7abf4.md" `
    "F2" "2eff61ea-3a49-4af9-a49a-115643b7abf4.06819624d60cdecd863e2321b8620db3234311b3.zh-cn.xlf" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/aa7f05aacdd104acef73026614c1040d384dc0b8/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/yuwzho/2eff61ea-3a49-4af9-a49a-115643b7abf4.06819624d60cdecd863e2321b8620db3234311b3.zh-cn.xlf" `
    "G2" "2016-01-07 11:07:29"

Set-HandbackRow $wsZh "B3" "E3" "78740cd6-903a-4fec-92db-091bdfe1fe15.md" `
    "https://github.com/OpenLocalizationTest/oltest/blob/7d4cba76265f0c0f490220592f21e3eec828d83f/e2e/78740cd6-903a-4fec-92db-091bdfe1fe15.md" `
    "F3" "78740cd6-903a-4fec-92db-091bdfe1fe15.4984dd8aa532fb358edde686fcdf65b9d5557d07.zh-cn.xlf" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/aa7f05aacdd104acef73026614c1040d384dc0b8/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/yuwzho/78740cd6-903a-4fec-92db-091bdfe1fe15.4984dd8aa532fb358edde686fcdf65b9d5557d07.zh-cn.xlf" `
    "G3" "2016-01-07 11:07:29"

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

Set-HandbackRow $wsDe "B2" "E2" "2eff61ea-3a49-4af9-a49a-115643b7abf4.md" `
    "https://github.com/OpenLocalizationTest/oltest/blob/7d4cba76265f0c0f490220592f21e3eec828d83f/e2e/2eff61ea-3a49-4af9-a49a-115643b7abf4.md" `
    "F2" "2eff61ea-3a49-4af9-a49a-115643b7abf4.06819624d60cdecd863e2321b8620db3234311b3.de-de.xlf" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/e3abcb479d2157ef2b473c09c4d2d0b0efb1da40/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/yuwzho/2eff61ea-3a49-4af9-a49a-115643b7abf4.06819624d60cdecd863e2321b8620db3234311b3.de-de.xlf" `
    "G2" "2016-01-07 11:07:11"

Set-HandbackRow $wsDe "B3" "E3" "78740cd6-903a-4fec-92db-091bdfe1fe15.md" `
    "https://github.com/OpenLocalizationTest/oltest/blob/7d4cba76265f0c0f490220592f21e3eec828d83f/e2e/78740cd6-903a-4fec-92db-091bdfe1fe15.md" `
    "F3" "78740cd6-903a-4fec-92db-091bdfe1fe15.4984dd8aa532fb358edde686fcdf65b9d5557d07.de-de.xlf" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/e3abcb479d2157ef2b473c09c4d2d0b0efb1da40/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/yuwzho/78740cd6-903a-4fec-92db-091bdfe1fe15.4984dd8aa532fb358edde686fcdf65b9d5557d07.de-de.xlf" `
    "G3" "2016-01-07 11:07:11"
